$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a date-valued cell (serial number) with the same number format
# used elsewhere in the sheet (m/d/yy -> builtin numFmtId 14), so the saved
# style matches the existing "s=1" cells exactly.
function Set-DateCell($row, $col, $serial) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = $serial
    $cell.NumberFormat = "m/d/yy"
}

function Set-Cell($row, $col, $val) {
    $ws.Cells.Item($row, $col).Value = $val
}

# ---- Row 494: PER record for person 1568 ----
Set-Cell 494 1 "PER"
Set-Cell 494 2 1568
Set-Cell 494 3 1960
Set-Cell 494 11 1
Set-Cell 494 12 1

# ---- Row 495: D_ERA record for person 1568 ----
Set-Cell 495 1 "D_ERA"
Set-Cell 495 2 1568
Set-Cell 495 4 1398937
Set-DateCell 495 5 39491
Set-DateCell 495 6 39565

# ---- Row 496: D_ERA record for person 1568 ----
Set-Cell 496 1 "D_ERA"
Set-Cell 496 2 1568
Set-Cell 496 4 902427
Set-DateCell 496 5 39491
Set-DateCell 496 6 39565

# ---- Row 497: VIS record for person 1568 ----
Set-Cell 497 1 "VIS"
Set-Cell 497 2 1568
Set-DateCell 497 5 39491
Set-DateCell 497 6 39565
Set-Cell 497 8 9201

# ---- Row 498: D_EXP record for person 1568 ----
Set-Cell 498 1 "D_EXP"
Set-Cell 498 2 1568
Set-Cell 498 4 40223504
Set-DateCell 498 5 39491
Set-DateCell 498 6 39565
Set-Cell 498 9 30
Set-Cell 498 10 30
Set-Cell 498 13 "2 times daily"
Set-Cell 498 14 "null"
Set-Cell 498 15 "null"

# ---- Row 499: D_EXP record for person 1568 ----
Set-Cell 499 1 "D_EXP"
Set-Cell 499 2 1568
Set-Cell 499 4 1594707
Set-DateCell 499 5 39491
Set-DateCell 499 6 39565
Set-Cell 499 9 30
Set-Cell 499 10 30
Set-Cell 499 13 "2 times daily"
Set-Cell 499 14 "null"
Set-Cell 499 15 "null"

# ---- Row 500: PER record for person 1569 ----
Set-Cell 500 1 "PER"
Set-Cell 500 2 1569
Set-Cell 500 3 1966
Set-Cell 500 11 1
Set-Cell 500 12 1

# ---- Row 501: D_ERA record for person 1569 ----
Set-Cell 501 1 "D_ERA"
Set-Cell 501 2 1569
Set-Cell 501 4 1398937
Set-DateCell 501 5 39491
Set-DateCell 501 6 39565

# ---- Row 502: D_ERA record for person 1569 ----
Set-Cell 502 1 "D_ERA"
Set-Cell 502 2 1569
Set-Cell 502 4 902427
Set-DateCell 502 5 39491
Set-DateCell 502 6 39565

# ---- Row 503: VIS record for person 1569 ----
Set-Cell 503 1 "VIS"
Set-Cell 503 2 1569
Set-DateCell 503 5 39491
Set-DateCell 503 6 39565
Set-Cell 503 8 9201

# ---- Row 504: D_EXP record for person 1569 ----
Set-Cell 504 1 "D_EXP"
Set-Cell 504 2 1569
Set-Cell 504 4 19079775
Set-DateCell 504 5 39491
Set-DateCell 504 6 39565
Set-Cell 504 9 20
Set-Cell 504 10 30
Set-Cell 504 13 "4 times daily"
Set-Cell 504 14 "null"
Set-Cell 504 15 "null"

# ---- Row 505: D_EXP record for person 1569 ----
Set-Cell 505 1 "D_EXP"
Set-Cell 505 2 1569
Set-Cell 505 4 902489
Set-DateCell 505 5 39491
Set-DateCell 505 6 39565
Set-Cell 505 9 5
Set-Cell 505 10 30
Set-Cell 505 13 "4 times daily"
Set-Cell 505 14 "null"
Set-Cell 505 15 "null"

# ---- Row 506: PER record for person 1570 ----
Set-Cell 506 1 "PER"
Set-Cell 506 2 1570
Set-Cell 506 3 1950
Set-Cell 506 11 1
Set-Cell 506 12 1

# ---- Row 507: D_ERA record for person 1570 ----
Set-Cell 507 1 "D_ERA"
Set-Cell 507 2 1570
Set-Cell 507 4 1398937
Set-DateCell 507 5 39491
Set-DateCell 507 6 39565

# ---- Row 508: D_ERA record for person 1570 ----
Set-Cell 508 1 "D_ERA"
Set-Cell 508 2 1570
Set-Cell 508 4 950370
Set-DateCell 508 5 39491
Set-DateCell 508 6 39565

# ---- Row 509: VIS record for person 1570 ----
Set-Cell 509 1 "VIS"
Set-Cell 509 2 1570
Set-DateCell 509 5 39491
Set-DateCell 509 6 39565
Set-Cell 509 8 9201

# ---- Row 510: D_EXP record for person 1570 ----
Set-Cell 510 1 "D_EXP"
Set-Cell 510 2 1570
Set-Cell 510 4 40223506
Set-DateCell 510 5 39491
Set-DateCell 510 6 39565
Set-Cell 510 9 10
Set-Cell 510 10 10
Set-Cell 510 13 "Daily"
Set-Cell 510 14 "null"
Set-Cell 510 15 "null"

# ---- Row 511: D_EXP record for person 1570 ----
Set-Cell 511 1 "D_EXP"
Set-Cell 511 2 1570
Set-Cell 511 4 43219718
Set-DateCell 511 5 39491
Set-DateCell 511 6 39565
Set-Cell 511 9 30
Set-Cell 511 10 30
Set-Cell 511 13 "2 times daily"
Set-Cell 511 14 "null"
Set-Cell 511 15 "null"

# Update selection / active cell to match the final view state.
$ws.Range("L507").Select()
